$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 289, shifting existing rows
# (289..343) down to (290..344).
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new weekly price record.
$ws.Range("A289").Value = 8
$ws.Range("B289").Value = "Terminal La Palmera de La Serena"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = 44694
$ws.Range("E289").Value = 4
$ws.Range("F289").Value = 100114013
$ws.Range("G289").Value = "Zanahoria"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 600
$ws.Range("K289").Value = 6000
$ws.Range("L289").Value = 7000
$ws.Range("M289").Value = 6500
$ws.Range("N289").Value = "$/saco 20 kilos"
$ws.Range("O289").Value = "Provincia del Elquí"
$ws.Range("P289").Value = 325
$ws.Range("Q289").Value = 20
$ws.Range("R289").Value = "Hortaliza"
